$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 244, shifting existing rows 244-267 down to 245-268
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with data
$ws.Cells.Item(244, 1).Value = 4
$ws.Cells.Item(244, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(244, 3).Value = "Los Lagos"
$ws.Cells.Item(244, 4).Value = 44918
$ws.Cells.Item(244, 5).Value = 10
$ws.Cells.Item(244, 6).Value = "Fruta"
$ws.Cells.Item(244, 7).Value = 100109
$ws.Cells.Item(244, 8).Value = "Uva"
$ws.Cells.Item(244, 9).Value = 100109001
$ws.Cells.Item(244, 10).Value = "Uva"
$ws.Cells.Item(244, 11).Value = "Superior Seedless"
$ws.Cells.Item(244, 12).Value = "Primera"
$ws.Cells.Item(244, 13).Value = 500
$ws.Cells.Item(244, 14).Value = 20000
$ws.Cells.Item(244, 15).Value = 21000
$ws.Cells.Item(244, 16).Value = 20500
$ws.Cells.Item(244, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(244, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(244, 19).Value = 2562
$ws.Cells.Item(244, 20).Value = 8
